$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 282; this shifts the existing rows
# 282..396 down to 283..397 (and carries the D-column date style
# down onto the new blank row).
$ws.Rows.Item(282).Insert()

# Populate the newly inserted row 282 with the new record.
$ws.Cells.Item(282, 1).Value = 10
$ws.Cells.Item(282, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(282, 3).Value = "La Araucanía"
$ws.Cells.Item(282, 4).Value = 45006
$ws.Cells.Item(282, 5).Value = 9
$ws.Cells.Item(282, 6).Value = "Fruta"
$ws.Cells.Item(282, 7).Value = 100102
$ws.Cells.Item(282, 8).Value = "Cítricos"
$ws.Cells.Item(282, 9).Value = 100102006
$ws.Cells.Item(282, 10).Value = "Pomelo"
$ws.Cells.Item(282, 11).Value = "Start Ruby"
$ws.Cells.Item(282, 12).Value = "Primera"
$ws.Cells.Item(282, 13).Value = 90
$ws.Cells.Item(282, 14).Value = 15000
$ws.Cells.Item(282, 15).Value = 15000
$ws.Cells.Item(282, 16).Value = 15000
$ws.Cells.Item(282, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(282, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(282, 19).Value = 1000
$ws.Cells.Item(282, 20).Value = 15
